# Change the presentation to use PowerPoint's current default 16:9
# ("widescreen") aspect ratio: 13.333in x 7.5in (12192000 x 6858000 EMU).
#
# Setting PageSetup.SlideSize to the "on-screen 16:9" constant is the
# COM-idiomatic equivalent of choosing Design > Slide Size > Widescreen
# in the UI; PowerPoint derives the corresponding SlideWidth/SlideHeight
# (960 x 540 points) from it automatically and tags the size as
# "screen16x9" rather than a bespoke custom size.
$p = $ppt.ActivePresentation

$ppSlideSizeOnScreen16x9 = 15
$p.PageSetup.SlideSize = $ppSlideSizeOnScreen16x9
